# DD_RFID_CM_v1.0.xlsx — "add UC alpha (1.3.1), chờ thầy xác nhận"
#
# 1. On sheet "Quy tắc đặt tên":
#      - update the naming-format cell (C8) to mention the new
#        [_alpha|_beta|_Realease] suffix
#      - insert a new example row (DD_RIFD_ScreenDesign_v1.0_beta) right
#        after the existing DD_RIFD_ScreenDesign_v1.0 example
# 2. Leave sheet "GIT" content untouched (its shared-string indices shift
#    automatically once a string is removed/added elsewhere — no edits
#    needed there).
# 3. The previously-active tab (GIT) becomes inactive; "Quy tắc đặt tên"
#    becomes the active tab/sheet, with its own last-used selection, while
#    GIT keeps a separate stored selection for when it's revisited.

$wb = $excel.ActiveWorkbook

$wsNaming = $wb.Worksheets.Item("Quy tắc đặt tên")
$wsGit    = $wb.Worksheets.Item("GIT")

# --- Update the format example in C8 ---------------------------------
$wsNaming.Range("C8").Value = "DD_RFID_<tên viết tắt tài liệu>_v<số phiên bản>[_alpha|_beta|_Realease]"

# --- Insert a new row for the extra "_beta" example -------------------
# Row15 already holds "DD_RIFD_ScreenDesign_v1.0"; row16 holds
# "DD_RIFD_URD_v3.2". Insert a blank row at 16 so the new beta example
# sits right under the ScreenDesign example, pushing the URD/SRS examples
# down by one row each.
$wsNaming.Rows("16:16").Insert() | Out-Null
$wsNaming.Range("C16").Value = "DD_RIFD_ScreenDesign_v1.0_beta"

# --- Update selections / active sheet ---------------------------------
# GIT keeps its own remembered selection (now K17) for whenever it is
# reselected, but it is no longer the active tab.
$wsGit.Activate() | Out-Null
$wsGit.Range("K17").Select() | Out-Null

# "Quy tắc đặt tên" becomes the active (last-activated) sheet/tab, with
# its own remembered selection (H13).
$wsNaming.Activate() | Out-Null
$wsNaming.Range("H13").Select() | Out-Null
